$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark (previously left after the
#    "Description" bullet in the Attributes list).
# ---------------------------------------------------------------------------
try {
    $old = $d.Bookmarks("_GoBack")
    $old.Delete()
} catch {
    # no pre-existing _GoBack bookmark - nothing to remove
}

# ---------------------------------------------------------------------------
# 2) "search engine" -> "search bar" (prose paragraph: "MainWindow is a
#    class that shows latest movies , recommended movies , users  and has
#    search engine functionality for movies.")
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$needle = "search engine"
$idx1 = $full.IndexOf($needle)
if ($idx1 -ge 0) {
    $r1 = $d.Range($idx1, $idx1 + $needle.Length)
    $r1.Text = "search bar"
}

# ---------------------------------------------------------------------------
# 3) "search engine" -> "search bar" (subtitle/heading run: "set_search_key
#    (key) -> sets search key for search engine")
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$idx2 = $full.LastIndexOf($needle)
if ($idx2 -ge 0) {
    $r2 = $d.Range($idx2, $idx2 + $needle.Length)
    $r2.Text = "search bar"
}

# ---------------------------------------------------------------------------
# 4) Re-add the "_GoBack" bookmark at the new edit location, right after
#    "... shows recommended movies " and before "if sort is specified ..."
#    inside the show_movies(**kwargs) description.
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$anchor = "shows recommended movies "
$idx3 = $full.IndexOf($anchor)
if ($idx3 -ge 0) {
    $pos = $idx3 + $anchor.Length
    $r3 = $d.Range($pos, $pos)
    $d.Bookmarks.Add("_GoBack", $r3)
}

Write-Output "edit complete"
